$d = $word.ActiveDocument

# Locate the paragraph containing the astromap link (old year 2018).
$oldText = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2018/)."
$newUrl  = "(http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

$rng = $d.Content
$found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, `
    $true, 1, $false, $null, 0)

if (-not $found) {
    throw "Could not find the astromap link text to update"
}

$para = $rng.Paragraphs(1)
$prng = $para.Range

# Pull the live paragraph-open-tag + <w:pPr> straight out of the document so
# we keep its exact attributes (paraId, rsids, ...) and formatting, and only
# replace the run content.
$oxml = $prng.WordOpenXML
if ($oxml -notmatch '(<w:p [^>]*>(<w:pPr>.*?</w:pPr>)?)') {
    throw "Could not read paragraph properties for the link paragraph"
}
$pOpenAndPPr = $matches[1]

$newParaXml = $pOpenAndPPr + '<w:r/><w:r><w:t>' + $newUrl + '</w:t></w:r></w:p>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' + `
    '<w:body>' + $newParaXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# Rebuild the paragraph with the original paragraph properties intact, but
# collapse the three old runs ("(", the hyperlink-styled URL, ").") into a
# single new run (with default/no explicit run formatting) holding the
# updated 2022 link, preceded by an empty run.
$prng.InsertXML($xml)

Write-Host "Updated astromap link year 2018 -> 2022"
